$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 13:22"

# Bizkaia/Vizcaya (row 7) - new counts, stays in rank position 7
$ws.Range("B7").Value = 5651
$ws.Range("C7").Value = 5092
$ws.Range("D7").Value = 4658
$ws.Range("E7").Value = 394

# Araba/Alava jumps above Toledo and Zaragoza into rank position 14
$ws.Range("A14").Value = "Araba/Alava"
$ws.Range("B14").Value = 2990
$ws.Range("C14").Value = 5092
$ws.Range("D14").Value = 4658
$ws.Range("E14").Value = 254

# Toledo shifts down to rank position 15 (its own figures unchanged)
$ws.Range("A15").Value = "Toledo"
$ws.Range("B15").Value = 2984
$ws.Range("C15").Value = 2205
$ws.Range("D15").Value = 9768
$ws.Range("E15").Value = 403

# Zaragoza shifts down to rank position 16 (its own figures unchanged)
$ws.Range("A16").Value = "Zaragoza"
$ws.Range("B16").Value = 2976
$ws.Range("C16").Value = 668
$ws.Range("D16").Value = 1993
$ws.Range("E16").Value = 315

# Gipuzkoa/Guipuzcoa jumps above Asturias into rank position 23
$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B23").Value = 1874
$ws.Range("C23").Value = 5092
$ws.Range("D23").Value = 4658
$ws.Range("E23").Value = 117

# Asturias shifts down to rank position 24 (its own figures unchanged)
$ws.Range("A24").Value = "Asturias"
$ws.Range("B24").Value = 1827
$ws.Range("C24").Value = 414
$ws.Range("D24").Value = 1285
$ws.Range("E24").Value = 128

# Caceres (row 28) - new counts, stays in rank position 28
$ws.Range("B28").Value = 1665
$ws.Range("C28").Value = 226
$ws.Range("D28").Value = 1190
$ws.Range("E28").Value = 249

# Badajoz (row 43) - new counts, stays in rank position 43
$ws.Range("B43").Value = 821
$ws.Range("C43").Value = 236
$ws.Range("D43").Value = 531
$ws.Range("E43").Value = 54

# Melilla (row 54) - new counts, stays in rank position 54
$ws.Range("B54").Value = 95
$ws.Range("C54").Value = 18
$ws.Range("D54").Value = 75
$ws.Range("E54").Value = 2
